# Apply updated values to existing rows (source data corrections)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95
$ws.Cells.Item(95, 3).Value = 20
$ws.Cells.Item(95, 5).Value = 24
$ws.Cells.Item(95, 6).Value = 816

# Row 96
$ws.Cells.Item(96, 3).Value = 19
$ws.Cells.Item(96, 5).Value = 21
$ws.Cells.Item(96, 6).Value = 886
$ws.Cells.Item(96, 7).Value = 1

# Row 102
$ws.Cells.Item(102, 3).Value = 49
$ws.Cells.Item(102, 6).Value = 1122
$ws.Cells.Item(102, 7).Value = 0

# Row 103
$ws.Cells.Item(103, 3).Value = 57
$ws.Cells.Item(103, 5).Value = 30
$ws.Cells.Item(103, 6).Value = 1103

# Row 109
$ws.Cells.Item(109, 3).Value = 67
$ws.Cells.Item(109, 5).Value = 26
$ws.Cells.Item(109, 6).Value = 1209

# Row 110
$ws.Cells.Item(110, 3).Value = 76
$ws.Cells.Item(110, 5).Value = 24
$ws.Cells.Item(110, 6).Value = 1171

# Row 116
$ws.Cells.Item(116, 3).Value = 119
$ws.Cells.Item(116, 5).Value = 21
$ws.Cells.Item(116, 6).Value = 1511
$ws.Cells.Item(116, 7).Value = 1

# Row 117
$ws.Cells.Item(117, 3).Value = 132
$ws.Cells.Item(117, 5).Value = 19
$ws.Cells.Item(117, 6).Value = 1503

# Row 123
$ws.Cells.Item(123, 3).Value = 197
$ws.Cells.Item(123, 5).Value = 41
$ws.Cells.Item(123, 6).Value = 1768

# Row 124
$ws.Cells.Item(124, 3).Value = 205
$ws.Cells.Item(124, 5).Value = 42
$ws.Cells.Item(124, 6).Value = 1837

# Row 130
$ws.Cells.Item(130, 5).Value = 61
$ws.Cells.Item(130, 6).Value = 2063

# Row 131
$ws.Cells.Item(131, 5).Value = 97
$ws.Cells.Item(131, 6).Value = 2054
$ws.Cells.Item(131, 7).Value = 3

# Row 133
$ws.Cells.Item(133, 3).Value = 251
$ws.Cells.Item(133, 5).Value = 100
$ws.Cells.Item(133, 6).Value = 1882
$ws.Cells.Item(133, 7).Value = 3

# Row 135
$ws.Cells.Item(135, 6).Value = 2264
$ws.Cells.Item(135, 7).Value = 6

# Row 137
$ws.Cells.Item(137, 5).Value = 123
$ws.Cells.Item(137, 6).Value = 2506
$ws.Cells.Item(137, 7).Value = 7

# Row 138
$ws.Cells.Item(138, 5).Value = 159
$ws.Cells.Item(138, 6).Value = 2807
$ws.Cells.Item(138, 7).Value = 8

# Row 139
$ws.Cells.Item(139, 5).Value = 156
$ws.Cells.Item(139, 6).Value = 2603

# Row 144
$ws.Cells.Item(144, 3).Value = 260
$ws.Cells.Item(144, 5).Value = 214
$ws.Cells.Item(144, 6).Value = 2962
$ws.Cells.Item(144, 7).Value = 5

# Row 145
$ws.Cells.Item(145, 3).Value = 285
$ws.Cells.Item(145, 5).Value = 219
$ws.Cells.Item(145, 6).Value = 2726

# Row 151
$ws.Cells.Item(151, 3).Value = 392
$ws.Cells.Item(151, 5).Value = 250
$ws.Cells.Item(151, 6).Value = 2279
$ws.Cells.Item(151, 7).Value = 15

# Row 152
$ws.Cells.Item(152, 5).Value = 280
$ws.Cells.Item(152, 6).Value = 2082

# Row 158
$ws.Cells.Item(158, 3).Value = 304
$ws.Cells.Item(158, 5).Value = 272
$ws.Cells.Item(158, 6).Value = 1952
$ws.Cells.Item(158, 7).Value = 5

# Row 159
$ws.Cells.Item(159, 3).Value = 278
$ws.Cells.Item(159, 5).Value = 282
$ws.Cells.Item(159, 6).Value = 1749
$ws.Cells.Item(159, 7).Value = 7

# Row 165
$ws.Cells.Item(165, 5).Value = 413
$ws.Cells.Item(165, 6).Value = 2249
$ws.Cells.Item(165, 7).Value = 5

# Row 166
$ws.Cells.Item(166, 3).Value = 369
$ws.Cells.Item(166, 5).Value = 394
$ws.Cells.Item(166, 6).Value = 2304
$ws.Cells.Item(166, 7).Value = 8

# Row 171
$ws.Cells.Item(171, 5).Value = 434
$ws.Cells.Item(171, 6).Value = 2382

# Row 172
$ws.Cells.Item(172, 5).Value = 386
$ws.Cells.Item(172, 6).Value = 2310
$ws.Cells.Item(172, 7).Value = 10

# Row 173
$ws.Cells.Item(173, 3).Value = 319
$ws.Cells.Item(173, 5).Value = 394
$ws.Cells.Item(173, 6).Value = 2326

# Row 187
$ws.Cells.Item(187, 5).Value = 301
$ws.Cells.Item(187, 6).Value = 1943
$ws.Cells.Item(187, 7).Value = 6

# Row 193
$ws.Cells.Item(193, 5).Value = 187
$ws.Cells.Item(193, 6).Value = 1449

# Row 194
$ws.Cells.Item(194, 5).Value = 207
$ws.Cells.Item(194, 6).Value = 1177

# Row 200
$ws.Cells.Item(200, 5).Value = 350
$ws.Cells.Item(200, 6).Value = 3395
$ws.Cells.Item(200, 7).Value = 11

# Row 201
$ws.Cells.Item(201, 5).Value = 338
$ws.Cells.Item(201, 6).Value = 3607
$ws.Cells.Item(201, 7).Value = 8

# Row 207
$ws.Cells.Item(207, 5).Value = 423
$ws.Cells.Item(207, 6).Value = 2866
$ws.Cells.Item(207, 7).Value = 12

# Row 208
$ws.Cells.Item(208, 5).Value = 428
$ws.Cells.Item(208, 6).Value = 2539
$ws.Cells.Item(208, 7).Value = 9

# Row 214
$ws.Cells.Item(214, 3).Value = 190
$ws.Cells.Item(214, 5).Value = 373
$ws.Cells.Item(214, 6).Value = 2595
$ws.Cells.Item(214, 7).Value = 9

# Row 215
$ws.Cells.Item(215, 5).Value = 320
$ws.Cells.Item(215, 6).Value = 2447
$ws.Cells.Item(215, 7).Value = 7

# Row 221
$ws.Cells.Item(221, 5).Value = 430
$ws.Cells.Item(221, 6).Value = 2556
$ws.Cells.Item(221, 7).Value = 5

# Row 222
$ws.Cells.Item(222, 5).Value = 426
$ws.Cells.Item(222, 6).Value = 2505
$ws.Cells.Item(222, 7).Value = 9

# Row 228
$ws.Cells.Item(228, 5).Value = 308
$ws.Cells.Item(228, 6).Value = 2876
$ws.Cells.Item(228, 7).Value = 12

# Row 229
$ws.Cells.Item(229, 5).Value = 272
$ws.Cells.Item(229, 6).Value = 2880
$ws.Cells.Item(229, 7).Value = 11

# Row 232
$ws.Cells.Item(232, 6).Value = 3895

# Row 292
$ws.Cells.Item(292, 3).Value = 752

# Row 303
$ws.Cells.Item(303, 3).Value = 431
$ws.Cells.Item(303, 5).Value = 581
$ws.Cells.Item(303, 6).Value = 2428
$ws.Cells.Item(303, 7).Value = 35

# Row 304
$ws.Cells.Item(304, 5).Value = 511
$ws.Cells.Item(304, 6).Value = 2411
$ws.Cells.Item(304, 7).Value = 37

# Row 305
$ws.Cells.Item(305, 3).Value = 417
$ws.Cells.Item(305, 5).Value = 536
$ws.Cells.Item(305, 6).Value = 2461
$ws.Cells.Item(305, 7).Value = 30

# Row 306
$ws.Cells.Item(306, 5).Value = 520
$ws.Cells.Item(306, 6).Value = 2574
$ws.Cells.Item(306, 7).Value = 36

# Row 307
$ws.Cells.Item(307, 3).Value = 616
$ws.Cells.Item(307, 5).Value = 546
$ws.Cells.Item(307, 6).Value = 2711

# Row 309
$ws.Cells.Item(309, 3).Value = 624
$ws.Cells.Item(309, 5).Value = 538
$ws.Cells.Item(309, 6).Value = 2517
$ws.Cells.Item(309, 7).Value = 34

# Row 310
$ws.Cells.Item(310, 3).Value = 688
$ws.Cells.Item(310, 5).Value = 737
$ws.Cells.Item(310, 6).Value = 3042
$ws.Cells.Item(310, 7).Value = 56

# Row 314
$ws.Cells.Item(314, 3).Value = 809

# Append new rows of data at the end (rows 315-317)
# Row 315
$cA = $ws.Cells.Item(315, 1)
$cA.Value = "'12.01.2021"
$cA.Style = "Normal"
$ws.Cells.Item(315, 2).Value = 116200
$ws.Cells.Item(315, 3).Value = 567
$ws.Cells.Item(315, 4).Value = 3280815
$ws.Cells.Item(315, 5).Value = 631
$ws.Cells.Item(315, 6).Value = 4386
$ws.Cells.Item(315, 7).Value = 15

# Row 316
$cA = $ws.Cells.Item(316, 1)
$cA.Value = "'13.01.2021"
$cA.Style = "Normal"
$ws.Cells.Item(316, 2).Value = 116668
$ws.Cells.Item(316, 3).Value = 468
$ws.Cells.Item(316, 4).Value = 3280815
$ws.Cells.Item(316, 5).Value = 692
$ws.Cells.Item(316, 6).Value = 4513
$ws.Cells.Item(316, 7).Value = 26

# Row 317
$cA = $ws.Cells.Item(317, 1)
$cA.Value = "'14.01.2021"
$cA.Style = "Normal"
$ws.Cells.Item(317, 2).Value = 117011
$ws.Cells.Item(317, 3).Value = 343
$ws.Cells.Item(317, 4).Value = 3280815
$ws.Cells.Item(317, 5).Value = 671
$ws.Cells.Item(317, 6).Value = 4782
$ws.Cells.Item(317, 7).Value = 21


